$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 969.70215
$ws.Range("I15").Value = 969.70215
$ws.Range("K15").Value = 2909.10645
$ws.Range("M15").Value = -2740.10645
$ws.Range("H19").Value = 1077.7142
$ws.Range("I19").Value = 539.8333
$ws.Range("J19").Value = 1481.125
$ws.Range("K19").Value = 539.8333
$ws.Range("L19").Value = 1481.125
$ws.Range("M19").Value = -364.8333
$ws.Range("N19").Value = -1831.125
$ws.Range("H32").Value = 598.125
$ws.Range("I32").Value = 401
$ws.Range("J32").Value = 626.2857
$ws.Range("K32").Value = 401
$ws.Range("L32").Value = 626.2857
$ws.Range("M32").Value = -75
$ws.Range("N32").Value = -1278.2857
$ws.Range("H40").Value = 1428.125
$ws.Range("I40").Value = 1353.2222
$ws.Range("J40").Value = 1524.4286
$ws.Range("K40").Value = 1353.2222
$ws.Range("L40").Value = 1524.4286
$ws.Range("M40").Value = -1178.2222
$ws.Range("N40").Value = -1874.4286
$ws.Range("H43").Value = 1985.1
$ws.Range("I43").Value = 5450.5
$ws.Range("K43").Value = 5450.5
$ws.Range("M43").Value = -5381.5
$ws.Range("H48").Value = 2950
$ws.Range("J48").Value = 2950
$ws.Range("L48").Value = 8850
$ws.Range("N48").Value = -9434
$ws.Range("H53").Value = 372.95
$ws.Range("I53").Value = 339
$ws.Range("K53").Value = 339
$ws.Range("M53").Value = 298
$ws.Range("H56").Value = 2950
$ws.Range("J56").Value = 2950
$ws.Range("L56").Value = 8850
$ws.Range("N56").Value = -9918
$ws.Range("H62").Value = 2221.0557
$ws.Range("I62").Value = 1781.8182
$ws.Range("J62").Value = 2911.2856
$ws.Range("K62").Value = 1781.8182
$ws.Range("L62").Value = 2911.2856
$ws.Range("M62").Value = -1157.8182
$ws.Range("N62").Value = -4159.2856
$ws.Range("H65").Value = 2221.0557
$ws.Range("I65").Value = 1781.8182
$ws.Range("J65").Value = 2911.2856
$ws.Range("K65").Value = 8909.091
$ws.Range("L65").Value = 14556.428
$ws.Range("M65").Value = -5789.091
$ws.Range("N65").Value = -20796.428
$ws.Range("H98").Value = 700
$ws.Range("I98").Value = 700
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 700
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = 798
$ws.Range("H100").Value = 1193.6957
$ws.Range("I100").Value = 1125.8235
$ws.Range("J100").Value = 1386
$ws.Range("K100").Value = 1125.8235
$ws.Range("L100").Value = 1386
$ws.Range("M100").Value = -584.8235
$ws.Range("N100").Value = -2468
$ws.Range("H107").Value = 867.5625
$ws.Range("I107").Value = 1020.3
$ws.Range("J107").Value = 613
$ws.Range("K107").Value = 1020.3
$ws.Range("L107").Value = 613
$ws.Range("M107").Value = 899.7
$ws.Range("N107").Value = -4453
$ws.Range("H111").Value = 5563934
$ws.Range("I111").Value = 10058.429
$ws.Range("K111").Value = 30175.287
$ws.Range("M111").Value = -27108.287
$ws.Range("H122").Value = 700
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2100
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 350
$ws.Range("H132").Value = 4390287
$ws.Range("I132").Value = 4721410.5
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 14164231.5
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -14161701.5
$ws.Range("N132").Value = -13760

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 72220.28999999999
$ws.Range("I2").Value = 883.875
$ws.Range("J2").Value = 167335.5
$ws.Range("K2").Value = 883.875
$ws.Range("L2").Value = 167335.5
$ws.Range("M2").Value = -770.875
$ws.Range("N2").Value = -167561.5
$ws.Range("H28").Value = 3988.2
$ws.Range("I28").Value = 3988.2
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 3988.2
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -3796.2
$ws.Range("H32").Value = 10444.35
$ws.Range("I32").Value = 9425.416999999999
$ws.Range("J32").Value = 19614.75
$ws.Range("K32").Value = 9425.416999999999
$ws.Range("L32").Value = 19614.75
$ws.Range("M32").Value = -9138.416999999999
$ws.Range("N32").Value = -20188.75
$ws.Range("H63").Value = 2331.875
$ws.Range("I63").Value = 1668.3334
$ws.Range("K63").Value = 1668.3334
$ws.Range("M63").Value = -982.3334
$ws.Range("H66").Value = 2331.875
$ws.Range("I66").Value = 1668.3334
$ws.Range("K66").Value = 8341.666999999999
$ws.Range("M66").Value = -4909.666999999999
$ws.Range("H99").Value = 3988.2
$ws.Range("I99").Value = 3988.2
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 3988.2
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -993.1999999999998
$ws.Range("H116").Value = 72220.28999999999
$ws.Range("I116").Value = 883.875
$ws.Range("J116").Value = 167335.5
$ws.Range("K116").Value = 883.875
$ws.Range("L116").Value = 167335.5
$ws.Range("M116").Value = 1410.125
$ws.Range("N116").Value = -171923.5
$ws.Range("H132").Value = 1910.907
$ws.Range("I132").Value = 1943.579
$ws.Range("J132").Value = 1662.6
$ws.Range("K132").Value = 5830.737
$ws.Range("L132").Value = 4987.799999999999
$ws.Range("M132").Value = -3300.737
$ws.Range("N132").Value = -10047.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 72220.28999999999
$ws.Range("I3").Value = 883.875
$ws.Range("J3").Value = 167335.5
$ws.Range("K3").Value = 883.875
$ws.Range("L3").Value = 167335.5
$ws.Range("M3").Value = -769.875
$ws.Range("N3").Value = -167563.5
$ws.Range("H134").Value = 1847.8125
$ws.Range("I134").Value = 1520.9762
$ws.Range("J134").Value = 4135.6665
$ws.Range("K134").Value = 4562.9286
$ws.Range("L134").Value = 12406.9995
$ws.Range("M134").Value = -2027.9286
$ws.Range("N134").Value = -17476.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 896.4167
$ws.Range("I16").Value = 846.8
$ws.Range("K16").Value = 846.8
$ws.Range("M16").Value = -559.8
$ws.Range("H22").Value = 404.1
$ws.Range("I22").Value = 417.625
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 417.625
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -67.625
$ws.Range("N22").Value = -1050
$ws.Range("H31").Value = 2512.2363
$ws.Range("I31").Value = 1724.6818
$ws.Range("J31").Value = 3037.2727
$ws.Range("K31").Value = 1724.6818
$ws.Range("L31").Value = 3037.2727
$ws.Range("M31").Value = -1429.6818
$ws.Range("N31").Value = -3627.2727
$ws.Range("H34").Value = 2512.2363
$ws.Range("I34").Value = 1724.6818
$ws.Range("J34").Value = 3037.2727
$ws.Range("K34").Value = 1724.6818
$ws.Range("L34").Value = 3037.2727
$ws.Range("M34").Value = -1522.6818
$ws.Range("N34").Value = -3441.2727
$ws.Range("H99").Value = 2739.739
$ws.Range("I99").Value = 2712
$ws.Range("J99").Value = 2770
$ws.Range("K99").Value = 2712
$ws.Range("L99").Value = 2770
$ws.Range("M99").Value = -1214
$ws.Range("N99").Value = -5766
$ws.Range("H113").Value = 896.4167
$ws.Range("I113").Value = 846.8
$ws.Range("K113").Value = 846.8
$ws.Range("M113").Value = 1323.2
$ws.Range("H126").Value = 2739.739
$ws.Range("I126").Value = 2712
$ws.Range("J126").Value = 2770
$ws.Range("K126").Value = 8136
$ws.Range("L126").Value = 8310
$ws.Range("M126").Value = -5666
$ws.Range("N126").Value = -13250
$ws.Range("H132").Value = 3837.4412
$ws.Range("I132").Value = 3544.9614
$ws.Range("J132").Value = 4788
$ws.Range("K132").Value = 10634.8842
$ws.Range("L132").Value = 14364
$ws.Range("M132").Value = -8104.8842
$ws.Range("N132").Value = -19424

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1813
$ws.Range("J39").Value = 1944.3
$ws.Range("L39").Value = 5832.9
$ws.Range("N39").Value = -6420.9
$ws.Range("H113").Value = 963.1667
$ws.Range("I113").Value = 2161.6667
$ws.Range("J113").Value = 563.6667
$ws.Range("K113").Value = 6485.000100000001
$ws.Range("L113").Value = 1691.0001
$ws.Range("M113").Value = -4315.000100000001
$ws.Range("N113").Value = -6031.0001
$ws.Range("H132").Value = 1748.7188
$ws.Range("I132").Value = 813
$ws.Range("J132").Value = 2060.625
$ws.Range("K132").Value = 7317
$ws.Range("L132").Value = 18545.625
$ws.Range("M132").Value = -4787
$ws.Range("N132").Value = -23605.625

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 19700
$ws.Range("J52").Value = 19700
$ws.Range("L52").Value = 19700
$ws.Range("N52").Value = -20218

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4474.75
$ws.Range("I22").Value = 9799.5
$ws.Range("J22").Value = 2699.8333
$ws.Range("K22").Value = 9799.5
$ws.Range("L22").Value = 2699.8333
$ws.Range("M22").Value = -9504.5
$ws.Range("N22").Value = -3289.8333
$ws.Range("H27").Value = 4474.75
$ws.Range("I27").Value = 9799.5
$ws.Range("J27").Value = 2699.8333
$ws.Range("K27").Value = 9799.5
$ws.Range("L27").Value = 2699.8333
$ws.Range("M27").Value = -9692.5
$ws.Range("N27").Value = -2913.8333
$ws.Range("H122").Value = 2689.4
$ws.Range("I122").Value = 2699.889
$ws.Range("K122").Value = 8099.667
$ws.Range("M122").Value = -5649.667

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2071.9092
$ws.Range("I122").Value = 1380.5385
$ws.Range("J122").Value = 3070.5557
$ws.Range("K122").Value = 4141.6155
$ws.Range("L122").Value = 9211.667099999999
$ws.Range("M122").Value = -1691.6155
$ws.Range("N122").Value = -14111.6671
$ws.Range("H136").Value = 1521.4865
$ws.Range("I136").Value = 553.2692
$ws.Range("K136").Value = 1659.8076
$ws.Range("M136").Value = 890.1924000000001
